# The source deck was produced by Aspose.Slides running in evaluation mode,
# which stamps every slide with a locked, non-selectable "Evaluation only." /
# "Created with Aspose.Slides for .NET ..." / "Copyright ..." watermark text
# box. This edit re-saves the presentation (e.g. after licensing / cleanup)
# with those watermark text boxes removed from every slide; all other shapes
# are left untouched.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shape = $slide.Shapes.Item($i)

        if ($shape.Name -eq "TextBox" -and $shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.StartsWith("Evaluation only.")) {
                $shape.Delete()
            }
        }
    }
}
